$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Shuttle Runs" column (F) entirely
$ws.Columns.Item(6).Delete()

# Rename "Chin ups" header to "Sprint 40y"
$ws.Range("E1").Value = "Sprint 40y"

# Remove the last data row (Julie Piers, row 4)
$ws.Rows.Item(4).Delete()
